$wb = $excel.ActiveWorkbook

# ===== Overview =====
$ws = $wb.Worksheets.Item("Overview")

# Update cell values
$ws.Range("A2").Value = "ffff92659b39-d24d-4ce8-bc4e-57332714bd88.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("A3").Value = "ffffff7c16e662-c219-47e4-9cac-459aba1abd40.md"
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("A4").Value = "569f2bbe-ada9-402d-892b-0544b541c87f.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("A5").Value = ".localization-config"
$ws.Range("B5").Value = "Not to be localized"
$ws.Range("C5").Value = "Not to be localized"

# Rebuild hyperlinks in the correct order so relationship ids realign
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/0d18e7bc545fa442db059e5db39233e9b56b317e/e2e/569f2bbe-ada9-402d-892b-0544b541c87f.md", [Type]::Missing, [Type]::Missing, "ffff92659b39-d24d-4ce8-bc4e-57332714bd88.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/0d18e7bc545fa442db059e5db39233e9b56b317e/e2e/ffff92659b39-d24d-4ce8-bc4e-57332714bd88.md", [Type]::Missing, [Type]::Missing, "ffffff7c16e662-c219-47e4-9cac-459aba1abd40.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/0d18e7bc545fa442db059e5db39233e9b56b317e/e2e/ffffff7c16e662-c219-47e4-9cac-459aba1abd40.md", [Type]::Missing, [Type]::Missing, "569f2bbe-ada9-402d-892b-0544b541c87f.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/0d18e7bc545fa442db059e5db39233e9b56b317e/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null

# ===== zh-cn =====
$ws = $wb.Worksheets.Item("zh-cn")

# Update cell values
$ws.Range("A2").Value = "ffff92659b39-d24d-4ce8-bc4e-57332714bd88.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.zh-cn.xlf"
$ws.Range("D2").Value = "2016-03-09 14:29:08"
$ws.Range("E2").Value = "3d6d923a-4f6f-4169-992f-ccc384019ff3.md"
$ws.Range("F2").Value = "3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.zh-cn.xlf"
$ws.Range("G2").Value = "2016-03-09 14:29:36"
$ws.Range("H2").Value = "Include"
$ws.Range("A3").Value = "ffffff7c16e662-c219-47e4-9cac-459aba1abd40.md"
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.zh-cn.xlf"
$ws.Range("D3").Value = "2016-03-09 14:29:08"
$ws.Range("E3").Value = "3d6d923a-4f6f-4169-992f-ccc384019ff3.md"
$ws.Range("F3").Value = "3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.zh-cn.xlf"
$ws.Range("G3").Value = "2016-03-09 14:29:36"
$ws.Range("H3").Value = "Include"
$ws.Range("A4").Value = "569f2bbe-ada9-402d-892b-0544b541c87f.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "569f2bbe-ada9-402d-892b-0544b541c87f.d6c4b3bb6e2795444b0c4976f08c6325dc8397f5.zh-cn.xlf"
$ws.Range("D4").Value = "2016-03-09 14:32:26"
$ws.Range("E4").Value = "569f2bbe-ada9-402d-892b-0544b541c87f.md"
$ws.Range("F4").Value = "569f2bbe-ada9-402d-892b-0544b541c87f.d6c4b3bb6e2795444b0c4976f08c6325dc8397f5.zh-cn.xlf"
$ws.Range("G4").Value = "2016-03-09 14:31:56"
$ws.Range("H4").Value = "Include"
$ws.Range("A5").Value = ".localization-config"
$ws.Range("B5").Value = "Not to be localized"
$ws.Range("D5").Value = "0001-01-01 00:00:00"
$ws.Range("G5").Value = "0001-01-01 00:00:00"
$ws.Range("H5").Value = "Ignored"

# Rebuild hyperlinks in the correct order so relationship ids realign
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/0d18e7bc545fa442db059e5db39233e9b56b317e/e2e/569f2bbe-ada9-402d-892b-0544b541c87f.md", [Type]::Missing, [Type]::Missing, "ffff92659b39-d24d-4ce8-bc4e-57332714bd88.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5dccd23a2698e2b85fe517cce745e1f95e7d20cc/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/569f2bbe-ada9-402d-892b-0544b541c87f.d6c4b3bb6e2795444b0c4976f08c6325dc8397f5.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/7450d45c935a8eda5802abda0f5251d5ae798451/e2e/569f2bbe-ada9-402d-892b-0544b541c87f.md", [Type]::Missing, [Type]::Missing, "3d6d923a-4f6f-4169-992f-ccc384019ff3.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/1aada4637e900a0020deeef4f3f29c39dbe0fd36/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/569f2bbe-ada9-402d-892b-0544b541c87f.d6c4b3bb6e2795444b0c4976f08c6325dc8397f5.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/0d18e7bc545fa442db059e5db39233e9b56b317e/e2e/ffff92659b39-d24d-4ce8-bc4e-57332714bd88.md", [Type]::Missing, [Type]::Missing, "ffffff7c16e662-c219-47e4-9cac-459aba1abd40.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8936cbf64460a8d88a1cbecad0156059d56b583f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/bfb945533f82a05b9749d4a053efda92fea8d709/e2e/3d6d923a-4f6f-4169-992f-ccc384019ff3.md", [Type]::Missing, [Type]::Missing, "3d6d923a-4f6f-4169-992f-ccc384019ff3.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/344e3e673f5f9c98f8f96d94a512be8902b61e14/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/0d18e7bc545fa442db059e5db39233e9b56b317e/e2e/ffffff7c16e662-c219-47e4-9cac-459aba1abd40.md", [Type]::Missing, [Type]::Missing, "569f2bbe-ada9-402d-892b-0544b541c87f.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8936cbf64460a8d88a1cbecad0156059d56b583f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "569f2bbe-ada9-402d-892b-0544b541c87f.d6c4b3bb6e2795444b0c4976f08c6325dc8397f5.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/bfb945533f82a05b9749d4a053efda92fea8d709/e2e/3d6d923a-4f6f-4169-992f-ccc384019ff3.md", [Type]::Missing, [Type]::Missing, "569f2bbe-ada9-402d-892b-0544b541c87f.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/344e3e673f5f9c98f8f96d94a512be8902b61e14/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "569f2bbe-ada9-402d-892b-0544b541c87f.d6c4b3bb6e2795444b0c4976f08c6325dc8397f5.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/0d18e7bc545fa442db059e5db39233e9b56b317e/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null

# ===== de-de =====
$ws = $wb.Worksheets.Item("de-de")

# Update cell values
$ws.Range("A2").Value = "ffff92659b39-d24d-4ce8-bc4e-57332714bd88.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.de-de.xlf"
$ws.Range("D2").Value = "2016-03-09 14:29:11"
$ws.Range("E2").Value = "3d6d923a-4f6f-4169-992f-ccc384019ff3.md"
$ws.Range("F2").Value = "3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.de-de.xlf"
$ws.Range("G2").Value = "2016-03-09 14:29:41"
$ws.Range("H2").Value = "Include"
$ws.Range("A3").Value = "ffffff7c16e662-c219-47e4-9cac-459aba1abd40.md"
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.de-de.xlf"
$ws.Range("D3").Value = "2016-03-09 14:29:11"
$ws.Range("E3").Value = "3d6d923a-4f6f-4169-992f-ccc384019ff3.md"
$ws.Range("F3").Value = "3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.de-de.xlf"
$ws.Range("G3").Value = "2016-03-09 14:29:41"
$ws.Range("H3").Value = "Include"
$ws.Range("A4").Value = "569f2bbe-ada9-402d-892b-0544b541c87f.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "569f2bbe-ada9-402d-892b-0544b541c87f.d6c4b3bb6e2795444b0c4976f08c6325dc8397f5.de-de.xlf"
$ws.Range("D4").Value = "2016-03-09 14:32:30"
$ws.Range("E4").Value = "569f2bbe-ada9-402d-892b-0544b541c87f.md"
$ws.Range("F4").Value = "569f2bbe-ada9-402d-892b-0544b541c87f.d6c4b3bb6e2795444b0c4976f08c6325dc8397f5.de-de.xlf"
$ws.Range("G4").Value = "2016-03-09 14:32:02"
$ws.Range("H4").Value = "Include"
$ws.Range("A5").Value = ".localization-config"
$ws.Range("B5").Value = "Not to be localized"
$ws.Range("D5").Value = "0001-01-01 00:00:00"
$ws.Range("G5").Value = "0001-01-01 00:00:00"
$ws.Range("H5").Value = "Ignored"

# Rebuild hyperlinks in the correct order so relationship ids realign
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/0d18e7bc545fa442db059e5db39233e9b56b317e/e2e/569f2bbe-ada9-402d-892b-0544b541c87f.md", [Type]::Missing, [Type]::Missing, "ffff92659b39-d24d-4ce8-bc4e-57332714bd88.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/903e1daa498e0975d6abd29c4c8d3a4bdb694168/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/569f2bbe-ada9-402d-892b-0544b541c87f.d6c4b3bb6e2795444b0c4976f08c6325dc8397f5.de-de.xlf", [Type]::Missing, [Type]::Missing, "3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/35ced1cf00ade5ea910fc4523be9794689330204/e2e/569f2bbe-ada9-402d-892b-0544b541c87f.md", [Type]::Missing, [Type]::Missing, "3d6d923a-4f6f-4169-992f-ccc384019ff3.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/d2b711e00559276affe93aff567935220aa7e3a5/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/569f2bbe-ada9-402d-892b-0544b541c87f.d6c4b3bb6e2795444b0c4976f08c6325dc8397f5.de-de.xlf", [Type]::Missing, [Type]::Missing, "3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/0d18e7bc545fa442db059e5db39233e9b56b317e/e2e/ffff92659b39-d24d-4ce8-bc4e-57332714bd88.md", [Type]::Missing, [Type]::Missing, "ffffff7c16e662-c219-47e4-9cac-459aba1abd40.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9294d13ee7f6855343eb677d18ff7b6b8dc09f84/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.de-de.xlf", [Type]::Missing, [Type]::Missing, "3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/a932d76745e8ce1a0dc1da24c37fe5b6aa9d96c8/e2e/3d6d923a-4f6f-4169-992f-ccc384019ff3.md", [Type]::Missing, [Type]::Missing, "3d6d923a-4f6f-4169-992f-ccc384019ff3.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9b206085512fda197301b56a2f586d7b0b06e0cb/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.de-de.xlf", [Type]::Missing, [Type]::Missing, "3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/0d18e7bc545fa442db059e5db39233e9b56b317e/e2e/ffffff7c16e662-c219-47e4-9cac-459aba1abd40.md", [Type]::Missing, [Type]::Missing, "569f2bbe-ada9-402d-892b-0544b541c87f.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9294d13ee7f6855343eb677d18ff7b6b8dc09f84/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.de-de.xlf", [Type]::Missing, [Type]::Missing, "569f2bbe-ada9-402d-892b-0544b541c87f.d6c4b3bb6e2795444b0c4976f08c6325dc8397f5.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/a932d76745e8ce1a0dc1da24c37fe5b6aa9d96c8/e2e/3d6d923a-4f6f-4169-992f-ccc384019ff3.md", [Type]::Missing, [Type]::Missing, "569f2bbe-ada9-402d-892b-0544b541c87f.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9b206085512fda197301b56a2f586d7b0b06e0cb/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/3d6d923a-4f6f-4169-992f-ccc384019ff3.3d05005b107a5dc333b298d708d845aa5536ffcc.de-de.xlf", [Type]::Missing, [Type]::Missing, "569f2bbe-ada9-402d-892b-0544b541c87f.d6c4b3bb6e2795444b0c4976f08c6325dc8397f5.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/0d18e7bc545fa442db059e5db39233e9b56b317e/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null

